$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.987.60'
$ws.Range("E2").Value = '  -1.59%  '

$ws.Range("D3").Value = '2.024.29'
$ws.Range("E3").Value = '  -2.71%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.31'
$ws.Range("E5").Value = '  -3.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.608'
$ws.Range("E6").Value = '  -4.55%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.99'
$ws.Range("E8").Value = '  -5.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.381'
$ws.Range("E9").Value = '  -2.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0792'
$ws.Range("E10").Value = '  +1.50%  '

$ws.Range("E11").Value = '  -3.63%  '

$ws.Range("D12").Value = '2.322.40'
$ws.Range("E12").Value = '  -2.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.28'
$ws.Range("E13").Value = '  -5.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.63'
$ws.Range("E14").Value = '  -2.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.745'
$ws.Range("E15").Value = '  -3.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.15'
$ws.Range("E16").Value = '  -3.92%  '

$ws.Range("D17").Value = '2.037.76'
$ws.Range("E17").Value = '  -2.15%  '

$ws.Range("D18").Value = '36.934.98'
$ws.Range("E18").Value = '  -1.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.03'
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.88'
$ws.Range("E20").Value = '  -2.65%  '

$ws.Range("D21").Value = '0.0₃0829'
$ws.Range("E21").Value = '  -0.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.76'
$ws.Range("E22").Value = '  -1.48%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("E24").Value = '  +3.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  -5.02%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.76'
$ws.Range("E26").Value = '  -1.67%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.33'
$ws.Range("E27").Value = '  -3.99%  '

$ws.Range("E28").Value = '  -4.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.79'
$ws.Range("E29").Value = '  -3.53%  '

$ws.Range("E30").Value = '  -2.93%  '

$ws.Range("E31").Value = '  -4.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.48'
$ws.Range("E32").Value = '  -3.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0610'
$ws.Range("E33").Value = '  -4.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.45'
$ws.Range("E34").Value = '  -4.17%  '

$ws.Range("E35").Value = '  -4.93%  '

$ws.Range("E36").Value = '  +0.34%  '

$ws.Range("E37").Value = '  +0.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.20'
$ws.Range("E38").Value = '  -4.31%  '

$ws.Range("E39").Value = '  +1.43%  '

$ws.Range("D40").Value = '1.500.97'
$ws.Range("E40").Value = '  +2.42%  '

$ws.Range("E41").Value = '  -5.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.92'
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '95.32'
$ws.Range("E43").Value = '  -5.32%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.82'
$ws.Range("E44").Value = '  -3.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0926'
$ws.Range("E45").Value = '  -3.36%  '

$ws.Range("E46").Value = '  -6.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.23'
$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -4.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("E49").Value = '  -1.45%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.211.54'
$ws.Range("E50").Value = '  -2.55%  '

$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.64'
$ws.Range("E51").Value = '  -9.53%  '
